# v1: save user setting for swatches-per-row
#
# The "program preference setting: how wide the palette area is ..." todo
# item (Id 20, row 3 of the "Active" sheet) has been implemented. Mark it
# Done by moving it from "Active" into the top data row of "Inactive"
# (pushing the other completed items down by one row).

$wb       = $excel.ActiveWorkbook
$active   = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Grab the pieces of the completed row we still need, using copy/paste so
# the original cell typing (text vs. shared string vs. number) is kept
# exactly as-is instead of being re-interpreted (e.g. date-like text
# being re-parsed into a date serial).

# Id, Title, Category carry straight over from "Active" row 3.
$active.Range("A3:B3").Copy()
$inactive.Rows.Item(2).Insert()
$inactive.Range("A2:B2").PasteSpecial()

$active.Range("D3").Copy()
$inactive.Range("D2").PasteSpecial()

# Created date (column E) also carries straight over.
$active.Range("E3").Copy()
$inactive.Range("E2").PasteSpecial()

# Status becomes "Done" (column C) and a "Done" date is recorded (column
# F) - both sourced from the row that used to be first so the same
# shared-string entries ("Done" / "8/11/2018") are reused.
$inactive.Range("C3").Copy()
$inactive.Range("C2").PasteSpecial()

$inactive.Range("F3").Copy()
$inactive.Range("F2").PasteSpecial()

# Drop the formatting Insert() copied down from the header row so the new
# row matches the plain data-row look of the rest of the sheet.
$inactive.Rows.Item(2).ClearFormats()

# Finally remove the now-duplicated row from "Active"; everything below
# it shifts up.
$active.Rows.Item(3).Delete()
